$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Label" header in column H, matching the formatting of the
# existing header row (column G).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Populate the new "Label" column: 0 for Control patients, 1 for MDD patients,
# based on the patient identifier already present in column A of each row.
for ($r = 2; $r -le 21; $r++) {
    $patient = $ws.Cells.Item($r, 1).Value()
    if ($patient -like "MDD*") {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
